# introducao-a-logica.pptx — apply the captured edit:
#   1. Every "datetimeFigureOut" date placeholder (on the slide master and
#      on all 11 slide layouts) is re-stamped from 08/04/2020 to 02/09/2021.
#   2. The title slide's title textbox drops the autofit font-scale back to
#      100% and the title text is shortened.

$p = $ppt.ActivePresentation

$oldDate = "08/04/2020"
$newDate = "02/09/2021"

# --- 1a. Slide master date placeholder -------------------------------------
$master = $p.SlideMaster
for ($j = 1; $j -le $master.Shapes.Count; $j++) {
    $sh = $master.Shapes.Item($j)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- 1b. Every slide layout's date placeholder ------------------------------
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $sh = $layout.Shapes.Item($j)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# --- 2. Title slide (slide 1) title shape -----------------------------------
$titleSlide = $p.Slides.Item(1)
$titleShape = $titleSlide.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Introdução ao raciocínio lógico"
# Drop the 90% autofit font-scale -> plain <a:normAutofit/>
$titleShape.TextFrame.AutoSize = 2
